# decomposition_main_te.xlsx — "tot_tut_validity & edits to the main paper"
#
# The sheet's data cells are formulas that pull cached values from an
# external workbook link ([1]decomposition_main_te!...). The headless
# engine has no access to that external source (and there is no COM
# surface to rewrite an external link's cached sheet data), so the
# updated figures are written directly onto the dependent cells as their
# new displayed text. A leading "'" (quote-prefix) is used wherever the
# new text would otherwise be auto-parsed as a number, so the stored
# cell keeps its original text semantics (no lost trailing zeros, no
# float round-off) - mirroring how `t="str"` cells looked before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5/6 (model row 5 of the source -> sheet rows 6/7) ---
$ws.Range("B6").Value = "-379.9***"
$ws.Range("E6").Value = "'-0.61"
$ws.Range("F6").Value = "-254.7**"
$ws.Range("G6").Value = "-0.064***"

$ws.Range("E7").Value = "(3.04)"
$ws.Range("F7").Value = "(104.7)"
$ws.Range("I7").Value = "(0.081)"

# --- Row 8/9 ---
$ws.Range("B8").Value = "'-78.8"
$ws.Range("E8").Value = "'-2.76"
$ws.Range("F8").Value = "'-55.3"
$ws.Range("G8").Value = "'-0.023"
$ws.Range("I8").Value = "'-0.10"

$ws.Range("B9").Value = "(114.5)"
$ws.Range("E9").Value = "(2.56)"
$ws.Range("F9").Value = "(109.1)"
$ws.Range("I9").Value = "(0.074)"

# --- Row 13 (Control mean) ---
$ws.Range("B13").Value = "'1850.6"
$ws.Range("E13").Value = "'5.75"
$ws.Range("F13").Value = "'1304.7"
$ws.Range("G13").Value = "'0.43"
$ws.Range("I13").Value = "'1.83"

# --- View/selection state ---
$ws.Range("A2:I13").Select()

# --- Row heights: the bordered header/footer rows shrink from 15.75 to 15
#     to match the new default row height used in this revision. ---
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
